$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade rows appended after the existing data (rows 7-9), matching the
# target diff: A=Principle, B=Start Principle, C=BuyPrice, D=SellPrice,
# E=IsShortSell, F=Price Change %, G=Date, H=Profitable.
$rows = @(
    @{ A = 9966.2800000000007; B = 10051.719999999999; C = 18.84; D = 19;    E = $true; F = 0.85; G = 42613.766944444447; H = $false },
    @{ A = 9945.35;            B = 9966.2800000000007;  C = 18.93; D = 18.97; E = $true; F = 0.21; G = 42614.675370370373; H = $false },
    @{ A = 9854.85;            B = 9945.35;             C = 18.72; D = 18.89; E = $true; F = 0.91; G = 42615.752118055556; H = $false }
)

$r = 7
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F

    # Copy the date-formatted cell above so the new cell inherits the same
    # style index (s="1") rather than minting a new numFmt/style entry.
    $ws.Cells.Item($r - 1, 7).Copy($ws.Cells.Item($r, 7))
    $ws.Cells.Item($r, 7).Value = $row.G

    $ws.Cells.Item($r, 8).Value = $row.H
    $r++
}
